$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "username"
$ws.Range("B3").Value = "lastname"
$ws.Range("C3").Value = "check1"

$ws.Range("D3").Select()
